$wb = $excel.ActiveWorkbook

$wsCap = $wb.Worksheets.Item("Capacitors")
$wsTests = $wb.Worksheets.Item("tests")

# --- Capacitors sheet ("Capacitors" = sheet1): add rows 5 and 6 ---
# Copy the date format (style) already used in column A of the tests sheet
# so the new date cells reuse the existing numFmtId 14 style instead of
# Excel manufacturing a brand-new custom number format.
$wsTests.Range("A3").Copy()
$wsCap.Range("A5").PasteSpecial(-4122)
$wsCap.Range("A5").Value = 41019

$wsTests.Range("A3").Copy()
$wsCap.Range("A6").PasteSpecial(-4122)
$wsCap.Range("A6").Value = 41019

$wsCap.Range("D5").Value = "0x0004"
$wsCap.Range("D6").Value = "0x0005"

$wsCap.Range("C5").Value = "Evens"
$wsCap.Range("C6").Value = "Evens"

$wsCap.Range("E5").Value = 4
$wsCap.Range("E6").Value = 5

$wsCap.Range("G5").Value = "100V"
$wsCap.Range("G6").Value = "300V"

# --- tests sheet (sheet2): add row 6 ---
$wsTests.Range("A3").Copy()
$wsTests.Range("A6").PasteSpecial(-4122)
$wsTests.Range("A6").Value = 41026

$wsTests.Range("A3").Copy()
$wsTests.Range("C6").PasteSpecial(-4122)
$wsTests.Range("C6").Value = 41026

$wsTests.Range("B6").Value = "buck filtering"
$wsTests.Range("D6").Value = "The capacitors were put in a low power buck supply to test for degredation"
$wsTests.Range("E6").Value = "0x0002"
$wsTests.Range("F6").Value = "0x0003"
$wsTests.Range("G6").Value = "0x0004"
$wsTests.Range("H6").Value = "0x0005"

# --- selections / active sheet ---
# Capacitors becomes the active/selected tab with cursor at G7;
# tests sheet keeps a (now-inactive) selection at A7.
$wsTests.Range("A7").Select() | Out-Null
$wsCap.Activate()
$wsCap.Range("G7").Select() | Out-Null
